$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row (row 2) with the new case/location/exposure info
$ws.Range("A2").Value = "Emerald"
$ws.Range("B2").Value = "Lakeside Paddle Boats, Emerald Lake Park"
$ws.Range("C2").Value = "31/12/20 3:30pm - 5:30pm"
$ws.Range("D2").Value = "Case visited venue"
$ws.Range("E2").Value = "old"

# Resize the columns to fit the new content (values chosen to match Excel's
# recalculated "best fit" widths for the new text)
$ws.Columns("A").ColumnWidth = 6.833333333333333
$ws.Columns("B").ColumnWidth = 33.0
$ws.Columns("C").ColumnWidth = 21.666666666666668
$ws.Columns("D").ColumnWidth = 14.333333333333334
$ws.Columns("E").ColumnWidth = 3.6666666666666665

# Select entire columns A:E, matching the saved selection state
$ws.Range("A1:E1048576").Select() | Out-Null

$wb.Save()
